# Reformat the Brinson attribution sheet:
#  - insert a new "Category" column (B) that classifies each row as
#    Industry / Currency / Country
#  - rename the old "Factors" header to "Segments" and relabel the three
#    group headers as "<TYPE> AGGREGATE"
#  - add four new "Country" segment rows (UK, USA, France, Germany) that
#    continue the numeric series with formulas
#  - add a trailing spacer row
#  - tidy up column widths / selection to match the new layout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header row (A1:F1)
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Segments"
$ws.Range("B1").Value = "Category"
$ws.Range("C1").Value = "Portfolio weight"
$ws.Range("D1").Value = "Benchmark weight"
$ws.Range("E1").Value = "Portfolio return"
$ws.Range("F1").Value = "Benchmark return"
$ws.Range("A1:F1").Font.Bold = $true

# ---------------------------------------------------------------------
# 2. Existing data rows (2-16): relabel column A, add column B category,
#    and shift the numeric series from B:E into C:F.
# ---------------------------------------------------------------------
$labels = @{
  2  = "INDUSTRY AGGREGATE"
  3  = "Benchmark"
  4  = "Healthcare"
  5  = "Technology"
  6  = "Media & Telecom"
  7  = "Cons. Disc"
  8  = "Cons. Staples"
  9  = "Energy"
  10 = "Industrials"
  11 = "CURRENCY AGGREGATE"
  12 = "USD"
  13 = "GBP"
  14 = "EUR"
  15 = "CAD"
  16 = "COUNTRY AGGREGATE"
}

$categories = @{
  2  = "Industry"; 3  = "Industry"; 4  = "Industry"; 5  = "Industry"
  6  = "Industry"; 7  = "Industry"; 8  = "Industry"; 9  = "Industry"; 10 = "Industry"
  11 = "Currency"; 12 = "Currency"; 13 = "Currency"; 14 = "Currency"; 15 = "Currency"
  16 = "Country"
}

$boldRows = @(2, 11, 16)

for ($r = 2; $r -le 16; $r++) {
  $ws.Range("A$r").Value = $labels[$r]
  $ws.Range("B$r").Value = $categories[$r]
  $ws.Range("C$r").Value = $r - 1
  $ws.Range("D$r").Value = $r - 1
  $ws.Range("E$r").Value = $r - 1
  $ws.Range("F$r").Value = $r - 1
}

foreach ($r in $boldRows) {
  $ws.Range("A$r").Font.Bold = $true
}

# ---------------------------------------------------------------------
# 3. New "Country" segment rows (17-20) with formulas continuing the
#    numeric sequence from row 16.
# ---------------------------------------------------------------------
$ws.Range("A17").Value = "UK"
$ws.Range("A18").Value = "USA"
$ws.Range("A19").Value = "France"
$ws.Range("A20").Value = "Germany"

$ws.Range("B17").Value = "Country"
$ws.Range("B18").Value = "Country"
$ws.Range("B19").Value = "Country"
$ws.Range("B20").Value = "Country"

$ws.Range("C17").Formula = "=C16+1"
$ws.Range("D17:F20").Formula = "=D16+1"
$ws.Range("C18:C20").Formula = "=C17+1"

# ---------------------------------------------------------------------
# 4. Trailing spacer row
# ---------------------------------------------------------------------
$ws.Range("C21").Value = " "

# ---------------------------------------------------------------------
# 5. Column widths / row outline levels / selection
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 21.85546875
$ws.Columns.Item(2).ColumnWidth = 15.7109375
$ws.Columns.Item(3).ColumnWidth = 15.5703125
$ws.Columns.Item(4).ColumnWidth = 17.7109375
$ws.Columns.Item(5).ColumnWidth = 15
$ws.Columns.Item(6).ColumnWidth = 17

$ws.Range("B13").Select()
